$d = $word.ActiveDocument

# --- Locate the "Pick up prompt ..." bullet paragraph and the two blank
#     paragraphs that follow it (end of the "Version 0.0.2" section,
#     right before the final sectPr). ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd() -eq "Pick up prompt no longer visible through other objects or walls") {
        $target = $i
        break
    }
}

$pPick = $d.Paragraphs($target)

# 1) Re-save the "Pick up prompt..." paragraph, splitting the trailing
#    word "walls" into its own run wrapped in a gramStart/gramEnd
#    proofing-error bookmark (matches Word's grammar-check pass over a
#    list item with no closing punctuation). The paragraph's own
#    w14:paraId/rsid attributes are carried over unchanged so only the
#    run content actually changes.
$pickXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
    'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" ' +
    'w14:paraId="7FC96979" w14:textId="38A433B9" w:rsidR="00FA18B1" w:rsidRPr="00D51198" ' +
    'w:rsidRDefault="00FA18B1" w:rsidP="00D51198">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Pick up prompt no longer visible through other objects or </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>walls</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '</w:p>'
$pPick.Range.InsertXML($pickXml)

# 2) The two following empty paragraphs collapse into a single new
#    bullet paragraph ("Texture Streaming Size to 2000Mb"). Delete the
#    first blank paragraph (merging it away) and rewrite the remaining
#    one (now the last paragraph in the body) with the new content.
$pBlank1 = $d.Paragraphs($target + 1)
$pBlank1.Range.Delete()

$pNew = $d.Paragraphs($target + 1)
$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Texture Streaming Size to </w:t></w:r>' +
    '<w:r><w:t>2000</w:t></w:r>' +
    '<w:r><w:t>Mb</w:t></w:r>' +
    '</w:p>'
$pNew.Range.InsertXML($newXml)

# 3) Bump the "Last Updated" footer timestamp field result.
$footer = $d.Sections(1).Footers(1)
$footer.Range.Find.Execute("19/01/2023 16:45", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "19/01/2023 23:28", 2)
